# "add Tim Booth ballot" — append a new voter row (row 25) to the ballots
# sheet, mirroring the layout of the existing rows: A=voter name, B..AJ are
# per-player "x" marks, AK=n_votes, AL=source, AM=date (mm/dd/yyyy cached as
# an Excel serial date).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 25

$ws.Cells.Item($row, 1).Value = "Tim Booth"   # A25 voter

# Columns Tim Booth voted for (marked with "x"), same encoding as every
# other row in the sheet: C=Barry Bonds, D=Roger Clemens, E=Roy Halladay,
# F=Todd Helton, I=Edgar Martinez, K=Mike Mussina, O=Mariano Rivera,
# P=Scott Rolen, Q=Curt Schilling, V=Larry Walker.
$votedCols = @(3, 4, 5, 6, 9, 11, 15, 16, 17, 22)
foreach ($col in $votedCols) {
    $ws.Cells.Item($row, $col).Value = "x"
}

$ws.Cells.Item($row, 37).Value = 10         # AK25 n_votes
$ws.Cells.Item($row, 38).Value = "Twitter"  # AL25 source
$ws.Cells.Item($row, 39).Value = 43441      # AM25 date (12/7/2018)

# Copy the number format from the date cell directly above (row 24) onto
# the new date cell so it reuses the sheet's existing date style instead of
# minting a brand-new custom numFmt/style entry.
$srcDate = $ws.Cells.Item($row - 1, 39)
$dstDate = $ws.Cells.Item($row, 39)
$srcDate.Copy()
$dstDate.PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Update the view: scroll so column M is the left-most visible column and
# select C25, matching the saved sheetView state.
$ws.Application.ActiveWindow.ScrollColumn = 13
$ws.Range("C25").Select()
